$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Gross Expenditures (D3) and Total Labor Cost (D5) values.
$ws.Range("D3").Value = 407429.16
$ws.Range("D5").Value = 67205.65

# Move the active selection to D4 (formulas D6, D7, D8, D9, E8, E9 recalc automatically).
$ws.Range("D4").Select()
